# Two-level deep TOC. Fixes after review with Aurimas.

$wb = $excel.ActiveWorkbook

# --- Add the new "SmartContractPrices" sheet at the end of the workbook ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws4 = $wb.Worksheets.Add($null, $lastSheet)
$ws4.Name = "SmartContractPrices"

# --- Column widths / formatting ---
$ws4.Columns.Item(1).ColumnWidth = 40.6640625
$ws4.Columns.Item(2).ColumnWidth = 26.6640625
$ws4.Columns.Item(3).ColumnWidth = 14.109375
$ws4.Columns.Item(4).ColumnWidth = 13.77734375
$ws4.Columns.Item(7).ColumnWidth = 10.33203125

# --- Header row (row 11) + data text, entered in the order the original
#     authoring session used (keeps the shared-string table layout aligned) ---
$ws4.Range("C11").Value = "Kuras"
$ws4.Range("D11").Value = "Kaina, ETH"
$ws4.Range("E11").Value = "Kaina, €"
$ws4.Range("B11").Value = "Vienetas"

$ws4.Range("B14").Value = "Funkcija getAttributeAsUser"
$ws4.Range("B13").Value = "Funkcija getAttribute"
$ws4.Range("B16").Value = "Funkcija removeAccess"
$ws4.Range("B15").Value = "Funkcija grantAccess"
$ws4.Range("B17").Value = "Funkcija requestAttributeAccess"
$ws4.Range("B12").Value = "Kontrakto sukūrimas"

$ws4.Range("A15").Value = 'Saugomas string "tZM11CdI7z4mZJc+/5kg3Q=="'
$ws4.Range("A16").Value = 'Kuris yra "+37063554865"'

$ws4.Range("G11").Value = 43232
$ws4.Range("G11").NumberFormat = "mm-dd-yy"

# --- Data rows (numbers) ---
$ws4.Range("C12").Value = 1074130
$ws4.Range("D12").Value = 0.0053707
$ws4.Range("E12").Value = 2.98611

$ws4.Range("C13").Value = 0
$ws4.Range("D13").Value = 0
$ws4.Range("E13").Value = 0

$ws4.Range("C14").Value = 0
$ws4.Range("D14").Value = 0
$ws4.Range("E14").Value = 0

$ws4.Range("C15").Value = 68499
$ws4.Range("D15").Value = 0.0003425
$ws4.Range("E15").Value = 0.19043

$ws4.Range("C16").Value = 51334
$ws4.Range("D16").Value = 0.0002567
$ws4.Range("E16").Value = 0.14273

$ws4.Range("C17").Value = 23034
$ws4.Range("D17").Value = 0.0001152
$ws4.Range("E17").Value = 0.06405

$ws4.Range("C18").Formula = "=SUM(C12:C17)"
$ws4.Range("D18:E18").Formula = "=SUM(D12:D17)"

# --- selection / view ---
$ws4.Range("M5").Select()
